$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 1 headers: add grammar (D1), synonyms (E1); shift definition_md to F1
$ws.Range("D1").Value = "grammar"
$ws.Range("E1").Value = "synonyms"
$ws.Range("F1").Value = "definition_md"

# Row 2
$ws.Range("A2").Value = "NCPED"
$ws.Range("B2").Value = "abbhantara"
$ws.Range("C2").Value = "interior, internal; being within, included in"
$ws.Range("D2").Value = "mfn. & neuter"
$ws.Range("E2").NumberFormat = "General"
$ws.Range("F2").Value = "1. (mfn.) interior, internal; being within, included in, among; belonging to one ‘s house, personal, intimate.`n2. (n.)`n   1. intermediate space, interval; the inside, interior.`n   2. a measure of length (= 28 hatthas)."

# Row 3
$ws.Range("A3").ClearContents()
$ws.Range("B3").Value = "ajjhokāse"
$ws.Range("C3").Value = "in the open air, in the open."
$ws.Range("D3").Value = "ind."
$ws.Range("E3").NumberFormat = "General"
$ws.Range("F3").Value = "in the open air, in the open."

# Row 4
$ws.Range("A4").ClearContents()
$ws.Range("B4").Value = "aṭṭita"
$ws.Range("C4").Value = "distressed; grieved; pained (see aṭṭiyati)"
$ws.Range("D4").Value = "pp mfn. "
$ws.Range("E4").NumberFormat = "General"
$ws.Range("F4").Value = "distressed; grieved; pained (see *[aṭṭiyati](/define/aṭṭiyati)*)"

# Row 5
$ws.Range("A5").Value = "NCPED"
$ws.Range("B5").Value = "cakkhuviññāṇa"
$ws.Range("C5").Value = "cognizance by the sense-organ that is the eye."
$ws.Range("D5").Value = "neuter"
$ws.Range("E5").NumberFormat = "General"
$ws.Range("F5").Value = "cognizance by the sense\-organ that is the eye."

# Row 6
$ws.Range("A6").Value = "NCPED"
$ws.Range("B6").Value = "cakkhuma"
$ws.Range("C6").Value = "possessing eyes, gifted with sight;"
$ws.Range("D6").Value = "mfn. & masculine"
$ws.Range("E6").NumberFormat = "General"
$ws.Range("F6").Value = "1. possessing eyes, gifted with sight; (one) who has eyes, who can see; one who has the gift of sight.`n2. one who possess insight and vision; wise."

# Row 7
$ws.Range("A7").Value = "NCPED"
$ws.Range("B7").Value = "ababa"
$ws.Range("C7").Value = "the name of a hell, or place in Avīci"
$ws.Range("D7").Value = "masculine"
$ws.Range("E7").NumberFormat = "General"
$ws.Range("F7").Value = "the name of a hell, or place in Avīci, where one suffers for an *ababa* of years."

# Wrap text for long definitions (matches previous wrap formatting, now on column F)
$ws.Range("F2").WrapText = $true
$ws.Range("F3").WrapText = $true
$ws.Range("F6").WrapText = $true

# Row heights
$ws.Rows.Item(2).RowHeight = 57.45
$ws.Rows.Item(3).RowHeight = 12.8
$ws.Rows.Item(4).RowHeight = 12.8
$ws.Rows.Item(5).RowHeight = 12.8
$ws.Rows.Item(6).RowHeight = 35.05
$ws.Rows.Item(7).RowHeight = 12.8

# Column widths (D=grammar, E=synonyms, F=definition_md)
$ws.Columns.Item(4).ColumnWidth = 14.75
$ws.Columns.Item(5).ColumnWidth = 8.92
$ws.Columns.Item(6).ColumnWidth = 50.75

# Selection moves to A4 per the saved view state
$ws.Range("A4").Select()

Write-Output "edit complete"
